$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select rows 7:10 (the Mouse/Rat experiment rows) and delete them entirely,
# shifting the Grand Offspring rows (old 11:16) up into rows 7:12.
$rng = $ws.Range("A7:XFD10")
$rng.Select()
$rng.EntireRow.Delete()
